$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A47").Value = "Giacomo Gasparini"
$ws.Range("B47").Value = "Mattia Kaiserman | Gli Introvabili"
$ws.Range("C47").Value = "Federico  Andreis | IMONTAGNA"
$ws.Range("D47").Value = "Michele Merighi | Clitoriders"
$ws.Range("E47").Value = "Geremia  Carollo | FC SAVIGNANO"
$ws.Range("F47").Value = "Alessio Delli Compagni | SdrumALA"
